# Auto-generated script applying the 2023-03-12 violent crime data update
# to the "violent-crime-full-year.xlsx" workbook.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("J2").Value = 1222
$ws.Range("J3").Value = 1306
$ws.Range("H4").Value = 1686
$ws.Range("I4").Value = 1757
$ws.Range("J4").Value = 289
$ws.Range("J6").Value = 1721
$ws.Range("H7").Value = 25999
$ws.Range("I7").Value = 26199
$ws.Range("J7").Value = 4632

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("J6").Value = 54
$ws.Range("J7").Value = 159

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("J2").Value = 34
$ws.Range("J6").Value = 57
$ws.Range("J7").Value = 172

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 31

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("J6").Value = 19
$ws.Range("J7").Value = 46

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("J8").Value = 289
$ws.Range("J9").Value = 30
$ws.Range("J11").Value = 57
$ws.Range("J13").Value = 7
$ws.Range("J15").Value = 53
$ws.Range("J18").Value = 63
$ws.Range("J19").Value = 162
$ws.Range("J25").Value = 25
$ws.Range("J29").Value = 263
$ws.Range("J31").Value = 31
$ws.Range("J33").Value = 192
$ws.Range("J34").Value = 30
$ws.Range("J36").Value = 66
$ws.Range("J37").Value = 159
$ws.Range("J39").Value = 4
$ws.Range("J42").Value = 189
$ws.Range("J44").Value = 38
$ws.Range("J47").Value = 40
$ws.Range("J48").Value = 32
$ws.Range("J52").Value = 103
$ws.Range("J54").Value = 92
$ws.Range("J60").Value = 29
$ws.Range("H63").Value = 237
$ws.Range("I63").Value = 188
$ws.Range("J63").Value = 17
$ws.Range("J64").Value = 27
$ws.Range("J67").Value = 172
$ws.Range("J69").Value = 12
$ws.Range("J71").Value = 22
$ws.Range("J72").Value = 19
$ws.Range("J76").Value = 77
$ws.Range("J77").Value = 36
$ws.Range("J78").Value = 63
$ws.Range("J79").Value = 137
$ws.Range("J84").Value = 46
$ws.Range("J85").Value = 196
$ws.Range("J88").Value = 36
$ws.Range("J90").Value = 51
$ws.Range("J92").Value = 15
$ws.Range("J94").Value = 36
$ws.Range("J98").Value = 33
$ws.Range("H101").Value = 25999
$ws.Range("I101").Value = 26199
$ws.Range("J101").Value = 4632

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("J3").Value = 53
$ws.Range("J6").Value = 79
$ws.Range("J7").Value = 192

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("J3").Value = 15
$ws.Range("J7").Value = 92

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("J2").Value = 74
$ws.Range("J3").Value = 99
$ws.Range("J6").Value = 69
$ws.Range("J7").Value = 263

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("J2").Value = 39
$ws.Range("J7").Value = 162

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Range("J4").Value = 5
$ws.Range("J7").Value = 38

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("J2").Value = 5
$ws.Range("J6").Value = 17
$ws.Range("J7").Value = 32

$ws = $wb.Worksheets.Item("River North")
$ws.Range("J6").Value = 45
$ws.Range("J7").Value = 77

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("J3").Value = 71
$ws.Range("J6").Value = 53
$ws.Range("J7").Value = 196

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("J2").Value = 35
$ws.Range("J7").Value = 189

$ws = $wb.Worksheets.Item("Boystown")
$ws.Range("J4").Value = 2
$ws.Range("J6").Value = 7

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("J4").Value = 9
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Norwood Park")
$ws.Range("J4").Value = 2
$ws.Range("J7").Value = 12

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("J3").Value = 48
$ws.Range("J6").Value = 39
$ws.Range("J7").Value = 137

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("J2").Value = 6
$ws.Range("J6").Value = 13
$ws.Range("J7").Value = 27

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("J6").Value = 38
$ws.Range("J7").Value = 63

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("J6").Value = 30
$ws.Range("J7").Value = 66

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("J2").Value = 25
$ws.Range("J3").Value = 32
$ws.Range("J7").Value = 103

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("J3").Value = 6
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("J6").Value = 21
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("J3").Value = 8
$ws.Range("J7").Value = 25

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("J2").Value = 9
$ws.Range("J7").Value = 40

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("J3").Value = 15
$ws.Range("J6").Value = 22
$ws.Range("J7").Value = 53

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("J3").Value = 4
$ws.Range("J7").Value = 33

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("J4").Value = 1
$ws.Range("J6").Value = 4

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("J2").Value = 18
$ws.Range("J7").Value = 57

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 30

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("J2").Value = 4
$ws.Range("J7").Value = 15

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("J3").Value = 13
$ws.Range("J7").Value = 36

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("J3").Value = 97
$ws.Range("J7").Value = 289

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("J3").Value = 14
$ws.Range("J7").Value = 51

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("J4").Value = 3
$ws.Range("J7").Value = 29

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("J6").Value = 14
$ws.Range("J7").Value = 22

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("J6").Value = 5
$ws.Range("J7").Value = 19

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("J3").Value = 11
$ws.Range("J7").Value = 36

